$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date serial for every data row (2-135).
# It was bumped by one day (2023-10-06 -> 2023-10-07, serial 45205 -> 45206).
for ($r = 2; $r -le 135; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}
